# Auto-generated edit script applying the OOXML diff changes
$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H5").Value = 2648.5
$ws.Range("I5").Value = 3356.3333
$ws.Range("J5").Value = 525
$ws.Range("K5").Value = 3356.3333
$ws.Range("L5").Value = 525
$ws.Range("M5").Value = -3241.3333
$ws.Range("N5").Value = -755
$ws.Range("H33").Value = 137
$ws.Range("I33").Value = 119.666664
$ws.Range("K33").Value = 119.666664
$ws.Range("M33").Value = 109.333336
$ws.Range("H55").Value = 445.875
$ws.Range("I55").Value = 100
$ws.Range("K55").Value = 100
$ws.Range("M55").Value = 114
$ws.Range("H127").Value = 1539.4
$ws.Range("J127").Value = 3000
$ws.Range("L127").Value = 9000
$ws.Range("N127").Value = -18920
$ws.Range("H137").Value = 3943
$ws.Range("I137").Value = 3879
$ws.Range("K137").Value = 11637
$ws.Range("M137").Value = -9087

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H102").Value = 1161.6428
$ws.Range("I102").Value = 989.46155
$ws.Range("K102").Value = 989.46155
$ws.Range("M102").Value = 632.53845
$ws.Range("H110").Value = 1150.5454
$ws.Range("J110").Value = 1000
$ws.Range("L110").Value = 1000
$ws.Range("N110").Value = -5090
$ws.Range("H122").Value = 4500
$ws.Range("I122").Value = 4500
$ws.Range("K122").Value = 13500
$ws.Range("M122").Value = -11050
$ws.Range("H132").Value = 1647.5
$ws.Range("I132").Value = 1647.5
$ws.Range("K132").Value = 4942.5
$ws.Range("M132").Value = -2412.5

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H76").Value = 13160
$ws.Range("J76").Value = 13160
$ws.Range("L76").Value = 13160
$ws.Range("N76").Value = -13790
$ws.Range("H79").Value = 13160
$ws.Range("J79").Value = 13160
$ws.Range("L79").Value = 13160
$ws.Range("N79").Value = -15344
$ws.Range("H105").Value = 3331.5715
$ws.Range("I105").Value = 3260.5557
$ws.Range("K105").Value = 3260.5557
$ws.Range("M105").Value = -1513.5557
$ws.Range("H107").Value = 1068.7646
$ws.Range("I107").Value = 869.5714
$ws.Range("K107").Value = 869.5714
$ws.Range("M107").Value = 1050.4286

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 3223.2666
$ws.Range("I31").Value = 2255.75
$ws.Range("J31").Value = 4329
$ws.Range("K31").Value = 2255.75
$ws.Range("L31").Value = 4329
$ws.Range("M31").Value = -1960.75
$ws.Range("N31").Value = -4919
$ws.Range("H34").Value = 3223.2666
$ws.Range("I34").Value = 2255.75
$ws.Range("J34").Value = 4329
$ws.Range("K34").Value = 2255.75
$ws.Range("L34").Value = 4329
$ws.Range("M34").Value = -2053.75
$ws.Range("N34").Value = -4733
$ws.Range("H105").Value = 1747.25
$ws.Range("I105").Value = 1495
$ws.Range("K105").Value = 1495
$ws.Range("M105").Value = 252
$ws.Range("H107").Value = 870.4666999999999
$ws.Range("I107").Value = 599.625
$ws.Range("K107").Value = 599.625
$ws.Range("M107").Value = 1320.375
$ws.Range("H122").Value = 900
$ws.Range("I122").Value = 900
$ws.Range("K122").Value = 2700
$ws.Range("M122").Value = -250
$ws.Range("H132").Value = 1296.8235
$ws.Range("I132").Value = 1003.13336
$ws.Range("J132").Value = 3499.5
$ws.Range("K132").Value = 3009.40008
$ws.Range("L132").Value = 10498.5
$ws.Range("M132").Value = -479.4000800000003
$ws.Range("N132").Value = -15558.5

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H7").Value = 318.5
$ws.Range("J7").Value = 425.66666
$ws.Range("L7").Value = 1276.99998
$ws.Range("N7").Value = -1500.99998
$ws.Range("H80").Value = 13779.7
$ws.Range("J80").Value = 14749.75
$ws.Range("L80").Value = 44249.25
$ws.Range("N80").Value = -46121.25
$ws.Range("H83").Value = 13779.7
$ws.Range("J83").Value = 14749.75
$ws.Range("L83").Value = 132747.75
$ws.Range("N83").Value = -142107.75
$ws.Range("H92").Value = 249.5
$ws.Range("I92").Value = 249.5
$ws.Range("K92").Value = 748.5
$ws.Range("M92").Value = 499.5

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H80").Value = 3597.6
$ws.Range("J80").Value = 3597.6
$ws.Range("L80").Value = 3597.6
$ws.Range("N80").Value = -5593.6
$ws.Range("H83").Value = 3597.6
$ws.Range("J83").Value = 3597.6
$ws.Range("L83").Value = 17988
$ws.Range("N83").Value = -27972
$ws.Range("H122").Value = 2123.9167
$ws.Range("I122").Value = 1347.6
$ws.Range("K122").Value = 4042.8
$ws.Range("M122").Value = -1592.8
$ws.Range("H126").Value = 18583
$ws.Range("I126").Value = 16299.6
$ws.Range("J126").Value = 30000
$ws.Range("K126").Value = 48898.8
$ws.Range("L126").Value = 90000
$ws.Range("M126").Value = -46428.8
$ws.Range("N126").Value = -94940
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H22").Value = 1500
$ws.Range("J22").Value = 1500
$ws.Range("L22").Value = 1500
$ws.Range("N22").Value = -2090
$ws.Range("H27").Value = 1500
$ws.Range("J27").Value = 1500
$ws.Range("L27").Value = 1500
$ws.Range("N27").Value = -1714
$ws.Range("H33").Value = 28800
$ws.Range("J33").Value = 28800
$ws.Range("L33").Value = 28800
$ws.Range("N33").Value = -29380

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H23").Value = 25190.25
$ws.Range("I23").Value = 33583.332
$ws.Range("K23").Value = 33583.332
$ws.Range("M23").Value = -33354.332
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H81").Value = 2773.75
$ws.Range("I81").Value = 2773.75
$ws.Range("K81").Value = 5547.5
$ws.Range("M81").Value = -4486.5
$ws.Range("H84").Value = 2773.75
$ws.Range("I84").Value = 2773.75
$ws.Range("K84").Value = 27737.5
$ws.Range("M84").Value = -22433.5
